$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update status column (H) for rows 22, 23, 25 ---
$ws.Range("H22").Value = "Completed"
$ws.Range("H23").Value = "In Progress"
$ws.Range("H25").Value = "In Progress"

# --- Add new task row 28 ---
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "Change comments on the order confirmation page and Order History Page"
$ws.Range("C28").Value = 0.5
$ws.Range("D28").Value = 41968
$ws.Range("D28").NumberFormat = "d-mmm"
$ws.Range("E28").Value = 41968
$ws.Range("E28").NumberFormat = "d-mmm"
$ws.Range("F28").Value = "Swarnima"
$ws.Range("G28").Value = "Add the comments - ""Kindly collect your order from Aisle 10"""

# --- Update view state ---
$ws.Range("G33").Select() | Out-Null
